# Update the "Percent Change" column (E) on the active sheet with the
# latest computed values. The sheet ships protected, so it must be
# unprotected before the cell values can be written.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("E2").Value = 0.01041028781383946
$ws.Range("E3").Value = 0.001375137513751179
$ws.Range("E4").Value = 0.009970674486803555
$ws.Range("E5").Value = -0.01625401625401623
$ws.Range("E6").Value = 0.001428611870793572
